$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# SamplesTab's TabQuery (cell B3) is rewritten: the query used to ORDER BY
# lab.test_name DESC; it now ORDER BYs smp.sample_id ASC instead.
# ---------------------------------------------------------------------------
$newSamplesQuery = @'
SELECT DISTINCT
    smp.sample_id AS "Sample ID",
    sts.study_subject_id AS "Case ID",
    prg.program_acronym AS "Program Code",
    std.study_acronym AS "ARM",
    sts.disease_subtype AS "Diagnosis",
    smp.tissue_type AS "Tissue Type",
    smp.composition AS "Tissue Composition",
    smp.sample_anatomic_site AS "Sample Anatomic Site",
    smp.method_of_sample_procurement AS "Sample Procurement Method",
    lab.test_name AS "platform"
FROM 
    df_program prg
LEFT JOIN 
    df_study std ON prg.program_id = std."program.program_id"
LEFT JOIN 
    df_study_subject sts ON std.study_id = sts."study.study_id"
LEFT JOIN
    df_diagnosis dgn ON sts.study_subject_id = dgn."study_subject.study_subject_id"
LEFT JOIN
    df_stratification_factor stf ON sts.study_subject_id = stf."study_subject.study_subject_id"
LEFT JOIN
    df_demographic_data dem ON sts.study_subject_id = dem."study_subject.study_subject_id"
LEFT JOIN
    df_sample smp ON sts.study_subject_id = smp."study_subject.study_subject_id"
LEFT JOIN
    df_laboratory_procedure lab ON prg.program_id = lab."program.program_id"
LEFT JOIN
    df_file fil ON smp.sample_id = fil."sample.sample_id"
LEFT JOIN
    df_therapeutic_procedure tpp ON dgn.diagnosis_id = tpp."diagnosis.diagnosis_id"
WHERE 
    tpp.chemotherapy_regimen = 'Dose dense AC (2 week cycles)'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
'@

$ws.Range("B3").Value2 = $newSamplesQuery

# Re-stamp B3's font (same visible size/color as the shared "wrap text" style
# used by the other query cells) so Excel records it as its own style entry,
# matching the distinct "s" index the real workbook ends up with after the
# cell was retyped.
$ws.Range("B3").Font.Size = 12
$ws.Range("B3").Font.ThemeColor = 1

# ---------------------------------------------------------------------------
# View state: scrolled down one row with C3 now the active cell, and the
# sheet's normal-view zoom reset to 100%.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C3").Select()
